{"js": "// Update the division-fact answers in the table to the new values.\n// Each [oldText, newText] pair corresponds to one <w:t> run in the table;\n// every old value is unique within the document, so a body-wide search\n// for the exact old text and replacing that single hit is safe.\nconst pairs = [\n  [\"57\u00f78=7, 1\", \"52\u00f78=6, 4\"],\n  [\"19\u00f78=2, 3\", \"73\u00f76=12, 1\"],\n  [\"86\u00f77=12, 2\", \"79\u00f73=26, 1\"],\n  [\"28\u00f75=5, 3\", \"40\u00f77=5, 5\"],\n  [\"37\u00f75=7, 2\", \"48\u00f73=16, 0\"],\n  [\"14\u00f79=1, 5\", \"23\u00f77=3, 2\"],\n  [\"37\u00f79=4, 1\", \"81\u00f78=10, 1\"],\n  [\"36\u00f72=18, 0\", \"46\u00f72=23, 0\"],\n  [\"30\u00f79=3, 3\", \"87\u00f78=10, 7\"],\n  [\"59\u00f74=14, 3\", \"13\u00f72=6, 1\"],\n  [\"52\u00f72=26, 0\", \"72\u00f77=10, 2\"],\n  [\"41\u00f78=5, 1\", \"45\u00f76=7, 3\"],\n  [\"31\u00f72=15, 1\", \"34\u00f73=11, 1\"],\n  [\"24\u00f79=2, 6\", \"48\u00f73=16, 0\"],\n  [\"87\u00f79=9, 6\", \"16\u00f75=3, 1\"],\n  [\"76\u00f74=19, 0\", \"28\u00f79=3, 1\"],\n  [\"76\u00f79=8, 4\", \"58\u00f77=8, 2\"],\n  [\"62\u00f74=15, 2\", \"30\u00f74=7, 2\"],\n  [\"55\u00f79=6, 1\", \"52\u00f75=10, 2\"],\n  [\"20\u00f78=2, 4\", \"78\u00f72=39, 0\"],\n  [\"52\u00f76=8, 4\", \"58\u00f79=6, 4\"],\n  [\"40\u00f75=8, 0\", \"69\u00f74=17, 1\"],\n  [\"88\u00f76=14, 4\", \"40\u00f77=5, 5\"],\n  [\"91\u00f76=15, 1\", \"79\u00f72=39, 1\"],\n  [\"21\u00f74=5, 1\", \"62\u00f78=7, 6\"],\n];\n\nconst body = context.document.body;\nconst allResults = [];\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  allResults.push({ results, newText });\n}\nawait context.sync();\n\nfor (const { results, newText } of allResults) {\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the division-fact answers in the table to the new values.\n# Each (old -> new) pair below corresponds to one <w:t> run in the table;\n# all old values are unique in the document, so a plain Find/Replace on the\n# whole document body is safe and will touch exactly one cell per pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('57\u00f78=7, 1', '52\u00f78=6, 4'),\n    @('19\u00f78=2, 3', '73\u00f76=12, 1'),\n    @('86\u00f77=12, 2', '79\u00f73=26, 1'),\n    @('28\u00f75=5, 3', '40\u00f77=5, 5'),\n    @('37\u00f75=7, 2', '48\u00f73=16, 0'),\n    @('14\u00f79=1, 5', '23\u00f77=3, 2'),\n    @('37\u00f79=4, 1', '81\u00f78=10, 1'),\n    @('36\u00f72=18, 0', '46\u00f72=23, 0'),\n    @('30\u00f79=3, 3', '87\u00f78=10, 7'),\n    @('59\u00f74=14, 3', '13\u00f72=6, 1'),\n    @('52\u00f72=26, 0', '72\u00f77=10, 2'),\n    @('41\u00f78=5, 1', '45\u00f76=7, 3'),\n    @('31\u00f72=15, 1', '34\u00f73=11, 1'),\n    @('24\u00f79=2, 6', '48\u00f73=16, 0'),\n    @('87\u00f79=9, 6', '16\u00f75=3, 1'),\n    @('76\u00f74=19, 0', '28\u00f79=3, 1'),\n    @('76\u00f79=8, 4', '58\u00f77=8, 2'),\n    @('62\u00f74=15, 2', '30\u00f74=7, 2'),\n    @('55\u00f79=6, 1', '52\u00f75=10, 2'),\n    @('20\u00f78=2, 4', '78\u00f72=39, 0'),\n    @('52\u00f76=8, 4', '58\u00f79=6, 4'),\n    @('40\u00f75=8, 0', '69\u00f74=17, 1'),\n    @('88\u00f76=14, 4', '40\u00f77=5, 5'),\n    @('91\u00f76=15, 1', '79\u00f72=39, 1'),\n    @('21\u00f74=5, 1', '62\u00f78=7, 6'),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
